# Apply the edits described by the commit:
#  - Add a new annotation row (row 23) on the "Non OCRED - OCR Problems" sheet:
#       A23 = 6.13 (page number, same number formatting/style as the rows above)
#       B23 = hyperlink-styled Issue Link URL
#       C23 = Comment "Article, contrast"
#    and extend the (till-now empty) style down to the blank row below (row 24).
#  - Make "Non OCRED - OCR Problems" the active/selected sheet tab, with A24 selected.
#  - "Pages Off" is consequently no longer the selected tab (its selection, C10,
#    is left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Non OCRED - OCR Problems")

# Copy the formatting of the last filled row in column A down onto the new
# row (23) and the following blank row (24), matching the right-aligned
# "page number" style used throughout column A.
$ws.Range("A22").Copy() | Out-Null
$ws.Range("A23:A24").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new annotation entry.
$ws.Range("C23").Value = "Article, contrast"
$ws.Range("B23").Value = "https://demo.humlab.umu.se/courier/081475engo.pdf"
$ws.Range("A23").Value = 6.13

# Make this sheet the active tab, with A24 as the selected cell.
$ws.Activate() | Out-Null
$ws.Range("A24").Select() | Out-Null
